$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Coluna" (column A) values: replace spaces/parentheses with underscores
# to match the new snake_case-style column naming convention.
$ws.Range("A2").Value = "Marital_status"
$ws.Range("A3").Value = "Application_mode"
$ws.Range("A4").Value = "Application_order"
$ws.Range("A6").Value = "Daytime/evening_attendance"
$ws.Range("A7").Value = "Previous_qualification"
$ws.Range("A8").Value = "Previous_qualification_grade"
$ws.Range("A10").Value = "Mother_qualification"
$ws.Range("A11").Value = "Father_qualification"
$ws.Range("A12").Value = "Mother_occupation"
$ws.Range("A13").Value = "Father_occupation"
$ws.Range("A14").Value = "Admission_grade"
$ws.Range("A16").Value = "Educational_special_needs"
$ws.Range("A18").Value = "Tuition_fees_up_to_date"
$ws.Range("A20").Value = "Scholarship_holder"
$ws.Range("A21").Value = "Age_at_enrollment"
$ws.Range("A23").Value = "Curricular_units_1st_sem_credited"
$ws.Range("A24").Value = "Curricular_units_1st_sem_enrolled"
$ws.Range("A25").Value = "Curricular_units_1st_sem_evaluations"
$ws.Range("A26").Value = "Curricular_units_1st_sem_approved"
$ws.Range("A27").Value = "Curricular_units_1st_sem_grade"
$ws.Range("A28").Value = "Curricular_units_1st_sem_without_evaluations"
$ws.Range("A29").Value = "Curricular_units_2nd_sem_credited"
$ws.Range("A30").Value = "Curricular_units_2nd_sem_enrolled"
$ws.Range("A31").Value = "Curricular_units_2nd_sem_evaluations"
$ws.Range("A32").Value = "Curricular_units_2nd_sem_approved"
$ws.Range("A33").Value = "Curricular_units_2nd_sem_grade"
$ws.Range("A34").Value = "Curricular_units_2nd_sem_without_evaluations"
$ws.Range("A35").Value = "Unemployment_rate"
$ws.Range("A36").Value = "Inflation_rate"
